$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append (dates 2021-09-02 through 2021-09-09)
$data = @(
    @(44441, 0, 6, 85.34850640113798),
    @(44442, 3, 7, 99.5732574679943),
    @(44443, 0, 7, 99.5732574679943),
    @(44444, 0, 6, 85.34850640113798),
    @(44445, 1, 7, 99.5732574679943),
    @(44446, 2, 7, 99.5732574679943),
    @(44447, 0, 6, 85.34850640113798),
    @(44448, 1, 7, 99.5732574679943)
)

$startRow = 367
$sourceStyleRange = $ws.Range("A366")

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Copy the formatting/style from the last existing data row (A366) to the new A-cell
    $sourceStyleRange.Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
